# Rapport du 22 Novembre 2025
# Append two new rows of data (row 29 and row 30) to the "Semaine_1" table
# on the "Semaine 1" worksheet, extend the table range, and update the
# sheet view to reflect the newly-entered rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Clone the formatting (font, number formats, wrap, etc.) of the last
#    existing data row (28) down into the two new rows (29-30) so the
#    new entries look like the rest of the table.
# ---------------------------------------------------------------------
$ws.Range("A28:P28").Copy()
$ws.Range("A29:P30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Write the new record's values into row 29 and row 30 (identical
#    content on both rows).
# ---------------------------------------------------------------------
foreach ($r in 29, 30) {
    $ws.Range("A$r").Value = 45983
    $ws.Range("B$r").Value = "Ndack NDAO"
    $ws.Range("C$r").Value = "GUEDIAWAYE"
    $ws.Range("D$r").Value = [char]0x0059 + "eumbeul Mb" + [char]0x00E9 + "de Sass"
    $ws.Range("E$r").Value = "MAMDOU DIA"
    $ws.Range("F$r").Value = 768059355
    $ws.Range("G$r").Value = "Grossiste"
    $ws.Range("H$r").Value = "Client Partenaire"
    $ws.Range("I$r").Value = "Livraison"
    $ws.Range("J$r").Value = "Ok"
    $ws.Range("K$r").Value = "Caf" + [char]0x00E9 + " stick Altimo 1,5gx09boites"
    $ws.Range("L$r").Value = 50
    $ws.Range("M$r").Value = 30000
    $ws.Range("N$r").Value = 1500000
    $ws.Range("O$r").Formula = "=""S""&_xlfn.ISOWEEKNUM(Semaine_1[[#This Row],[Date]])"
    $ws.Range("P$r").Formula = "=TEXT(Semaine_1[[#This Row],[Date]],""MMMM"")"
}

# ---------------------------------------------------------------------
# 3) Extend the "Semaine_1" table (and its AutoFilter) to cover the two
#    new rows.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item("Semaine_1")
$lo.Resize($ws.Range("A1:P30"))

# ---------------------------------------------------------------------
# 3b) The calculated "Semaine" / "Mois" columns on the pre-existing rows
#     lose their (invisible) background fill flag once the table is
#     refreshed with the new rows -- clear it explicitly to match.
# ---------------------------------------------------------------------
$ws.Range("O2:O28").Interior.Pattern = -4142
$ws.Range("P2:P28").Interior.Pattern = -4142

# ---------------------------------------------------------------------
# 4) Update the sheet view: select the freshly added rows, matching
#    where the user ended up after entering the new data.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A29:N30").Select()
